$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the computed totals in row 5 (D5, E5)
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 4

# Set up the sheet for printing (commit message: "Update to print out")
$ws.PageSetup.Orientation = 1

# Move the active selection to E6 (matches saved cursor position)
[void]$ws.Range("E6").Select()
